# Updated symbol list on Thu Jan 19 12:51:22 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for the crypto table.
# Values are written as text (matching the sheet's existing inline-string
# cells), so NumberFormat is forced to "@" before the write and the format
# is cleared again afterwards to avoid leaving a stray style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '288.96' },
    @{ Cell = 'E2'; Value = '-4.21%' },
    @{ Cell = 'D3'; Value = '30.70' },
    @{ Cell = 'E3'; Value = '-4.88%' },
    @{ Cell = 'D4'; Value = '4.938' },
    @{ Cell = 'E4'; Value = '-1.20%' },
    @{ Cell = 'D5'; Value = '0.07172' },
    @{ Cell = 'E5'; Value = '-6.48%' },
    @{ Cell = 'D6'; Value = '1.817' },
    @{ Cell = 'E6'; Value = '-10.86%' },
    @{ Cell = 'D7'; Value = '7.616' },
    @{ Cell = 'E7'; Value = '-2.92%' },
    @{ Cell = 'D8'; Value = '3.732' },
    @{ Cell = 'E8'; Value = '-1.34%' },
    @{ Cell = 'D9'; Value = '0.8968' },
    @{ Cell = 'E9'; Value = '-2.07%' },
    @{ Cell = 'D10'; Value = '0.1659' },
    @{ Cell = 'E10'; Value = '-5.69%' },
    @{ Cell = 'D11'; Value = '0.07731' },
    @{ Cell = 'E11'; Value = '-2.14%' },
    @{ Cell = 'D12'; Value = '0.07942' },
    @{ Cell = 'E12'; Value = '-6.39%' },
    @{ Cell = 'D13'; Value = '0.03038' },
    @{ Cell = 'E13'; Value = '-1.12%' },
    @{ Cell = 'D14'; Value = '0.1001' },
    @{ Cell = 'E14'; Value = '0.28%' },
    @{ Cell = 'D15'; Value = '0.001497' },
    @{ Cell = 'E15'; Value = '-0.94%' },
    @{ Cell = 'D16'; Value = '0.005782' },
    @{ Cell = 'E16'; Value = '0.94%' },
    @{ Cell = 'D18'; Value = '3.465' },
    @{ Cell = 'E18'; Value = '-0.05%' },
    @{ Cell = 'D19'; Value = '2.074' },
    @{ Cell = 'E19'; Value = '-3.65%' },
    @{ Cell = 'D20'; Value = '0.3318' },
    @{ Cell = 'E20'; Value = '-0.65%' },
    @{ Cell = 'D21'; Value = '0.1282' },
    @{ Cell = 'E21'; Value = '-3.50%' },
    @{ Cell = 'D22'; Value = '4.025' },
    @{ Cell = 'E22'; Value = '-5.58%' },
    @{ Cell = 'D23'; Value = '0.2103' },
    @{ Cell = 'E23'; Value = '5.68%' },
    @{ Cell = 'D24'; Value = '0.04514' },
    @{ Cell = 'E24'; Value = '-0.42%' },
    @{ Cell = 'D25'; Value = '0.001216' },
    @{ Cell = 'E25'; Value = '-1.12%' },
    @{ Cell = 'D26'; Value = '0.004617' },
    @{ Cell = 'E26'; Value = '4.83%' },
    @{ Cell = 'D27'; Value = '0.0001302' },
    @{ Cell = 'E27'; Value = '4.21%' },
    @{ Cell = 'D39'; Value = '0.01562' },
    @{ Cell = 'E39'; Value = '-8.58%' },
    @{ Cell = 'D40'; Value = '0.04339' },
    @{ Cell = 'E40'; Value = '-7.16%' },
    @{ Cell = 'D41'; Value = '0.007306' },
    @{ Cell = 'E41'; Value = '-3.27%' },
    @{ Cell = 'D42'; Value = '0.01005' },
    @{ Cell = 'D43'; Value = '0.1300' },
    @{ Cell = 'E43'; Value = '-3.80%' },
    @{ Cell = 'D44'; Value = '0.002063' },
    @{ Cell = 'E44'; Value = '-11.41%' },
    @{ Cell = 'D45'; Value = '0.009139' },
    @{ Cell = 'E45'; Value = '-13.15%' },
    @{ Cell = 'D46'; Value = '0.00005976' },
    @{ Cell = 'E46'; Value = '-4.43%' },
    @{ Cell = 'D47'; Value = '0.00000000751' },
    @{ Cell = 'E47'; Value = '0.17%' },
    @{ Cell = 'D48'; Value = '2.255' },
    @{ Cell = 'E48'; Value = '174.84%' },
    @{ Cell = 'E49'; Value = '0.01%' },
    @{ Cell = 'D50'; Value = '0.00002103' },
    @{ Cell = 'E50'; Value = '0.17%' },
    @{ Cell = 'D51'; Value = '0.0002003' },
    @{ Cell = 'E51'; Value = '0.17%' }

)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}
